$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.93
$ws.Range("N10").Value = 10
$ws.Range("V22").Value = 1.67
$ws.Range("U23").Value = 1.77
$ws.Range("V23").Value = 1.92
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = 2.1
$ws.Range("L24").Value = 3.75
$ws.Range("N24").Value = 9
$ws.Range("O24").Value = 1.33
$ws.Range("S24").Value = 1.44
$ws.Range("T24").Value = 2.63
$ws.Range("U24").Value = 1.77
$ws.Range("V24").Value = 1.92
$ws.Range("W24").Value = 7.5
$ws.Range("X24").Value = 11
$ws.Range("Z24").Value = 21
$ws.Range("AA24").Value = 19
$ws.Range("AB24").Value = 29
$ws.Range("AC24").Value = 9
$ws.Range("AE24").Value = 15
$ws.Range("AF24").Value = 51
$ws.Range("AG24").Value = 301
$ws.Range("AH24").Value = 9
$ws.Range("AI24").Value = 15
$ws.Range("AJ24").Value = 11
$ws.Range("AK24").Value = 34
$ws.Range("AL24").Value = 26
$ws.Range("AM24").Value = 34
$ws.Range("AN24").Value = 4.33
$ws.Range("AO24").Value = 13
$ws.Range("AP24").Value = 23
$ws.Range("AQ24").Value = 41
$ws.Range("AR24").Value = 67
$ws.Range("AS24").Value = 151
$ws.Range("AT24").Value = 2.63
$ws.Range("AU24").Value = 8
$ws.Range("AV24").Value = 51
$ws.Range("AW24").Value = 5
$ws.Range("AY24").Value = 26
$ws.Range("AZ24").Value = 51
$ws.Range("BA24").Value = 81
$ws.Range("BB24").Value = 201
$ws.Range("U25").Value = 1.63
$ws.Range("K26").Value = 2.05
$ws.Range("L26").Value = 3.75
$ws.Range("M26").Value = 1.07
$ws.Range("N26").Value = 8.5
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
$ws.Range("S26").Value = 1.44
$ws.Range("T26").Value = 2.63
$ws.Range("U26").Value = 1.83
$ws.Range("V26").Value = 1.83
$ws.Range("W26").Value = 7.5
$ws.Range("X26").Value = 11
$ws.Range("Y26").Value = 9.5
$ws.Range("AC26").Value = 8.5
$ws.Range("AE26").Value = 15
$ws.Range("AG26").Value = 301
$ws.Range("AH26").Value = 9
$ws.Range("AL26").Value = 26
$ws.Range("AP26").Value = 23
$ws.Range("AT26").Value = 2.63
$ws.Range("AU26").Value = 8
$ws.Range("AV26").Value = 51
$ws.Range("AX26").Value = 17
$ws.Range("AZ26").Value = 51
$ws.Range("BA26").Value = 81
$ws.Range("BB26").Value = 201
$ws.Range("G35").Value = 4.15
$ws.Range("H35").Value = 3.1
$ws.Range("I35").Value = 1.9
$ws.Range("J35").Value = 4.65
$ws.Range("K35").Value = 1.93
$ws.Range("N35").Value = 7.4
$ws.Range("P35").Value = 2.35
$ws.Range("Q35").Value = 2.32
$ws.Range("W35").Value = 9
$ws.Range("X35").Value = 21
$ws.Range("AB35").Value = 65
$ws.Range("AC35").Value = 6.7
$ws.Range("AD35").Value = 6.2
$ws.Range("AE35").Value = 19.5
$ws.Range("AH35").Value = 5.3
$ws.Range("AI35").Value = 7.6
$ws.Range("AK35").Value = 16
$ws.Range("AN35").Value = 5.7
$ws.Range("AQ35").Value = 150
$ws.Range("AU35").Value = 8.25
$ws.Range("AW35").Value = 3.55
$ws.Range("AX35").Value = 10
$ws.Range("AY35").Value = 24
$ws.Range("BA35").Value = 100
$ws.Range("M36").Value = 1.01
$ws.Range("O36").Value = 1.1
$ws.Range("M37").Value = 1.02
$ws.Range("O37").Value = 1.13
$ws.Range("M38").Value = 1.02
$ws.Range("O38").Value = 1.13
$ws.Range("P38").Value = 5
$ws.Range("M39").Value = 1.02
$ws.Range("O39").Value = 1.13
$ws.Range("M40").Value = 1.07
$ws.Range("N40").Value = 8.5
$ws.Range("AB40").Value = 34
$ws.Range("BB40").Value = 351
$ws.Range("J41").Value = 2.2
$ws.Range("M41").Value = 1.06
$ws.Range("N41").Value = 10
$ws.Range("AE41").Value = 15
$ws.Range("AJ41").Value = 19
$ws.Range("AS41").Value = 126
$ws.Range("M42").Value = 1.05
$ws.Range("N42").Value = 11
$ws.Range("W42").Value = 6
$ws.Range("AA42").Value = 13
$ws.Range("AB42").Value = 41
$ws.Range("AC42").Value = 11
$ws.Range("AE42").Value = 34
$ws.Range("AM42").Value = 81
$ws.Range("AN42").Value = 3
$ws.Range("AS42").Value = 201
$ws.Range("BA42").Value = 351
$ws.Range("K43").Value = 2.5
$ws.Range("M43").Value = 1.05
$ws.Range("N43").Value = 11
$ws.Range("W43").Value = 6
$ws.Range("X43").Value = 5
$ws.Range("AQ43").Value = 15
